$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2629091505494614
$ws.Range("C2").Value = 0.04063551494367346
$ws.Range("E2").Value = 0.1643923017644084
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002432897258490295
$ws.Range("I2").Value = 0.5344660747400489
$ws.Range("K2").Value = 0.2694279784494711
$ws.Range("M2").Value = 0.2198893434283633
$ws.Range("O2").Value = 2.340736204142956
$ws.Range("B3").Value = 0.2305834671728917
$ws.Range("C3").Value = 0.03627817298286118
$ws.Range("E3").Value = 0.1535769587299285
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.00243508909735331
$ws.Range("I3").Value = 0.542137487809546
$ws.Range("K3").Value = 0.2352035900325262
$ws.Range("M3").Value = 0.1975289479377267
$ws.Range("O3").Value = 2.368469090824803
$ws.Range("B4").Value = 0.2106926681933601
$ws.Range("C4").Value = 0.03358632794785876
$ws.Range("E4").Value = 0.147056273282125
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002436505613548645
$ws.Range("I4").Value = 0.5471805920492336
$ws.Range("K4").Value = 0.2141237625049115
$ws.Range("M4").Value = 0.1838598537537592
$ws.Range("O4").Value = 2.387049353816806
$ws.Range("B5").Value = 0.2025767713170978
$ws.Range("C5").Value = 0.03248529472270434
$ws.Range("E5").Value = 0.1444289999082571
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002437100689860821
$ws.Range("I5").Value = 0.5493193112742425
$ws.Range("K5").Value = 0.2055174777476623
$ws.Range("M5").Value = 0.178304760040561
$ws.Range("O5").Value = 2.395010832798178
$ws.Range("B6").Value = 0.2012285294987635
$ws.Range("C6").Value = 0.0323022236126036
$ws.Range("E6").Value = 0.1439945464538468
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002437200580590105
$ws.Range("I6").Value = 0.5496794917131727
$ws.Range("K6").Value = 0.2040874537986923
$ws.Range("M6").Value = 0.1773832578866674
$ws.Range("O6").Value = 2.396356364913501
$ws.Range("B7").Value = 0.2105832551425522
$ws.Range("C7").Value = 0.03357149550821248
$ws.Range("E7").Value = 0.1470207199304667
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.00243651356643737
$ws.Range("I7").Value = 0.5472090971222237
$ws.Range("K7").Value = 0.2140077596095153
$ws.Range("M7").Value = 0.1837848743596453
$ws.Range("O7").Value = 2.387155147110434
$ws.Range("B8").Value = 0.2517724077812602
$ws.Range("C8").Value = 0.03913653650302251
$ws.Range("E8").Value = 0.1606381395221916
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.002433638358848593
$ws.Range("I8").Value = 0.5370420535246456
$ws.Range("K8").Value = 0.2576413734365133
$ws.Range("M8").Value = 0.2121669702208848
$ws.Range("O8").Value = 2.349976094174565
$ws.Range("B9").Value = 0.3321872778399211
$ws.Range("C9").Value = 0.04991784908125396
$ws.Range("E9").Value = 0.1883048242875631
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002428558786458334
$ws.Range("I9").Value = 0.5197486755265928
$ws.Range("K9").Value = 0.3426662984888367
$ws.Range("M9").Value = 0.2683047487026968
$ws.Range("O9").Value = 2.289403261523972
$ws.Range("B10").Value = 0.3910318872014784
$ws.Range("C10").Value = 0.05775740686748065
$ws.Range("E10").Value = 0.209235927322176
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002425164035265847
$ws.Range("I10").Value = 0.5086595308193687
$ws.Range("K10").Value = 0.4047872291576198
$ws.Range("M10").Value = 0.3098496156781891
$ws.Range("O10").Value = 2.252447271022618
$ws.Range("B11").Value = 0.4177469747199893
$ws.Range("C11").Value = 0.06130593543119289
$ws.Range("E11").Value = 0.2188931175877045
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002423692175796249
$ws.Range("I11").Value = 0.5039666865079226
$ws.Range("K11").Value = 0.4329690376270605
$ws.Range("M11").Value = 0.3288165552036446
$ws.Range("O11").Value = 2.237279619666296
$ws.Range("B12").Value = 0.4278551471108756
$ws.Range("C12").Value = 0.0626470885409276
$ws.Range("E12").Value = 0.2225697958639969
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002423145181836403
$ws.Range("I12").Value = 0.502240284588737
$ws.Range("K12").Value = 0.4436292278356575
$ws.Range("M12").Value = 0.3360086701056844
$ws.Range("O12").Value = 2.231772900159328
$ws.Range("B13").Value = 0.425678546898439
$ws.Range("C13").Value = 0.06235836352441027
$ws.Range("E13").Value = 0.2217770772853527
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.00242326252625862
$ws.Range("I13").Value = 0.5026098405718997
$ws.Range("K13").Value = 0.4413338916639304
$ws.Range("M13").Value = 0.334459286435802
$ws.Range("O13").Value = 2.232948323773911
$ws.Range("B14").Value = 0.4185787486066488
$ws.Range("C14").Value = 0.06141632524146701
$ws.Range("E14").Value = 0.2191952038384528
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002423646966712342
$ws.Range("I14").Value = 0.503823638136442
$ws.Range("K14").Value = 0.4338462943891841
$ws.Range("M14").Value = 0.3294080596474203
$ws.Range("O14").Value = 2.236821826671132
$ws.Range("B15").Value = 0.4142288265584284
$ws.Range("C15").Value = 0.06083895986398602
$ws.Range("E15").Value = 0.2176163051832702
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002423883796140553
$ws.Range("I15").Value = 0.5045737269467594
$ws.Range("K15").Value = 0.429258392087462
$ws.Range("M15").Value = 0.3263153058592252
$ws.Range("O15").Value = 2.239225332419039
$ws.Range("B16").Value = 0.3892848681247756
$ws.Range("C16").Value = 0.05752514132200304
$ws.Range("E16").Value = 0.2086075532913299
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002425261678265245
$ws.Range("I16").Value = 0.5089733047947362
$ws.Range("K16").Value = 0.4029438804523977
$ws.Range("M16").Value = 0.3086114469612156
$ws.Range("O16").Value = 2.253471638541427
$ws.Range("B17").Value = 0.3739684383785402
$ws.Range("C17").Value = 0.055487646245723
$ws.Range("E17").Value = 0.20311584491958
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002426125480900366
$ws.Range("I17").Value = 0.5117624492656603
$ws.Range("K17").Value = 0.3867805989045223
$ws.Range("M17").Value = 0.2977680922761508
$ws.Range("O17").Value = 2.262632706486116
$ws.Range("B18").Value = 0.3651538182490697
$ws.Range("C18").Value = 0.05431406743834089
$ws.Range("E18").Value = 0.1999699026883732
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002426629137628922
$ws.Range("I18").Value = 0.5133997851645056
$ws.Range("K18").Value = 0.377476661243719
$ws.Range("M18").Value = 0.2915376808683163
$ws.Range("O18").Value = 2.268056600427087
$ws.Range("B19").Value = 0.3621684945498487
$ws.Range("C19").Value = 0.0539164292555796
$ws.Range("E19").Value = 0.1989069226086215
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002426800840157547
$ws.Range("I19").Value = 0.513959839575854
$ws.Range("K19").Value = 0.3743252796070067
$ws.Range("M19").Value = 0.2894292705678154
$ws.Range("O19").Value = 2.269919588731156
$ws.Range("B20").Value = 0.3755994214156146
$ws.Range("C20").Value = 0.05570471380427477
$ws.Range("E20").Value = 0.2036991261777814
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002426032822197343
$ws.Range("I20").Value = 0.5114621142138525
$ws.Range("K20").Value = 0.3885019615700003
$ws.Range("M20").Value = 0.2989217237808006
$ws.Range("O20").Value = 2.261641481253477
$ws.Range("B21").Value = 0.4206643597070752
$ws.Range("C21").Value = 0.06169309560797842
$ws.Range("E21").Value = 0.2199530259701703
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002423533766230623
$ws.Range("I21").Value = 0.5034657401687426
$ws.Range("K21").Value = 0.4360459035827375
$ws.Range("M21").Value = 0.3308914630996682
$ws.Range("O21").Value = 2.235677649987892
$ws.Range("B22").Value = 0.4500685570073983
$ws.Range("C22").Value = 0.06559168589024011
$ws.Range("E22").Value = 0.2306908762250259
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002421960895503779
$ws.Range("I22").Value = 0.4985350781300006
$ws.Range("K22").Value = 0.4670504739672197
$ws.Range("M22").Value = 0.3518423334737832
$ws.Range("O22").Value = 2.220090118732315
$ws.Range("B23").Value = 0.4343795851337973
$ws.Range("C23").Value = 0.06351233833963477
$ws.Range("E23").Value = 0.2249492889998095
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.00242279485476333
$ws.Range("I23").Value = 0.5011395977630428
$ws.Range("K23").Value = 0.4505091612932404
$ws.Range("M23").Value = 0.3406552717464493
$ws.Range("O23").Value = 2.228282899743292
$ws.Range("B24").Value = 0.3748620819806376
$ws.Range("C24").Value = 0.05560658439770805
$ws.Range("E24").Value = 0.2034353894981322
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002426074691265586
$ws.Range("I24").Value = 0.5115977903519635
$ws.Range("K24").Value = 0.3877237691409334
$ws.Range("M24").Value = 0.2984001558346847
$ws.Range("O24").Value = 2.262089124942875
$ws.Range("B25").Value = 0.3104730703174425
$ws.Range("C25").Value = 0.04701544239696887
$ws.Range("E25").Value = 0.1807153683529563
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002429873484279421
$ws.Range("I25").Value = 0.5241434676949126
$ws.Range("K25").Value = 0.3197244732690478
$ws.Range("M25").Value = 0.2530657724327057
$ws.Range("O25").Value = 2.304466459349513
